$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = 111271296
$ws.Range("Q2").Value2 = 473220.1559155915
$ws.Range("R2").Value2 = 6863539.25170773
$ws.Range("A3").Value2 = 111271923
$ws.Range("B3").Value2 = 96348
$ws.Range("D3").Value2 = "VU"
$ws.Range("E3").Value2 = 220787
$ws.Range("F3").Value2 = "Knärot"
$ws.Range("G3").Value2 = "Goodyera repens"
$ws.Range("H3").Value2 = "(L.) R. Br."
$ws.Range("Q3").Value2 = 473118.5439814709
$ws.Range("R3").Value2 = 6863582.939962601
$ws.Range("AC3").Value2 = "Tre blommande."
$ws.Range("A4").Value2 = 111270892
$ws.Range("B4").Value2 = 78578
$ws.Range("D4").Value2 = "NT"
$ws.Range("E4").Value2 = 6458
$ws.Range("F4").Value2 = "Lunglav"
$ws.Range("G4").Value2 = "Lobaria pulmonaria"
$ws.Range("H4").Value2 = "(L.) Hoffm."
$ws.Range("Q4").Value2 = 473239.9383552746
$ws.Range("R4").Value2 = 6863714.420922431
$ws.Range("A5").Value2 = 111270939
$ws.Range("B5").Value2 = 78578
$ws.Range("E5").Value2 = 6458
$ws.Range("F5").Value2 = "Lunglav"
$ws.Range("G5").Value2 = "Lobaria pulmonaria"
$ws.Range("H5").Value2 = "(L.) Hoffm."
$ws.Range("A6").Value2 = 111271723
$ws.Range("B6").Value2 = 96348
$ws.Range("D6").Value2 = "VU"
$ws.Range("E6").Value2 = 220787
$ws.Range("F6").Value2 = "Knärot"
$ws.Range("G6").Value2 = "Goodyera repens"
$ws.Range("H6").Value2 = "(L.) R. Br."
$ws.Range("Q6").Value2 = 473140.3516782348
$ws.Range("R6").Value2 = 6863595.022241795
$ws.Range("A7").Value2 = 111271309
$ws.Range("B7").Value2 = 78579
$ws.Range("D7").Value2 = "NT"
$ws.Range("E7").Value2 = 2081
$ws.Range("F7").Value2 = "Skrovellav"
$ws.Range("G7").Value2 = "Lobaria scrobiculata"
$ws.Range("H7").Value2 = "(Scop.) DC."
$ws.Range("Q7").Value2 = 473221.4734201821
$ws.Range("R7").Value2 = 6863586.84377678
$ws.Range("A8").Value2 = 111271068
$ws.Range("Q8").Value2 = 473238.8676645419
$ws.Range("R8").Value2 = 6863638.079474191
$ws.Range("A10").Value2 = 111270747
$ws.Range("Q10").Value2 = 473194.7999623233
$ws.Range("R10").Value2 = 6863736.454484907
$ws.Range("AC10").Value2 = "Sex blommande."
$ws.Range("A11").Value2 = 111271029
$ws.Range("B11").Value2 = 78579
$ws.Range("E11").Value2 = 2081
$ws.Range("F11").Value2 = "Skrovellav"
$ws.Range("G11").Value2 = "Lobaria scrobiculata"
$ws.Range("H11").Value2 = "(Scop.) DC."
$ws.Range("Q11").Value2 = 473229.5908188519
$ws.Range("R11").Value2 = 6863658.889402787
$ws.Range("A12").Value2 = 111270596
$ws.Range("B12").Value2 = 96348
$ws.Range("D12").Value2 = "VU"
$ws.Range("E12").Value2 = 220787
$ws.Range("F12").Value2 = "Knärot"
$ws.Range("G12").Value2 = "Goodyera repens"
$ws.Range("H12").Value2 = "(L.) R. Br."
$ws.Range("Q12").Value2 = 473184.8241620373
$ws.Range("R12").Value2 = 6863788.37406126
$ws.Range("AC12").Value2 = "Fem blommande."
$ws.Range("A13").Value2 = 111272062
$ws.Range("B13").Value2 = 96348
$ws.Range("D13").Value2 = "VU"
$ws.Range("E13").Value2 = 220787
$ws.Range("F13").Value2 = "Knärot"
$ws.Range("G13").Value2 = "Goodyera repens"
$ws.Range("H13").Value2 = "(L.) R. Br."
$ws.Range("Q13").Value2 = 473156.3705774212
$ws.Range("R13").Value2 = 6863531.269191674
$ws.Range("A14").Value2 = 111272375
$ws.Range("B14").Value2 = 96251
$ws.Range("D14").Value2 = "LC"
$ws.Range("E14").Value2 = 220093
$ws.Range("F14").Value2 = "Korallrot"
$ws.Range("G14").Value2 = "Corallorhiza trifida"
$ws.Range("H14").Value2 = "Châtel."
$ws.Range("Q14").Value2 = 473400.7315261344
$ws.Range("R14").Value2 = 6863573.187783281
$ws.Range("A16").Value2 = 111272292
$ws.Range("B16").Value2 = 78578
$ws.Range("E16").Value2 = 6458
$ws.Range("F16").Value2 = "Lunglav"
$ws.Range("G16").Value2 = "Lobaria pulmonaria"
$ws.Range("H16").Value2 = "(L.) Hoffm."
$ws.Range("Q16").Value2 = 473321.1690919191
$ws.Range("R16").Value2 = 6863539.403128584
$ws.Range("A17").Value2 = 111272343
$ws.Range("B17").Value2 = 78578
$ws.Range("D17").Value2 = "NT"
$ws.Range("E17").Value2 = 6458
$ws.Range("F17").Value2 = "Lunglav"
$ws.Range("G17").Value2 = "Lobaria pulmonaria"
$ws.Range("H17").Value2 = "(L.) Hoffm."
$ws.Range("Q17").Value2 = 473387.8703240218
$ws.Range("R17").Value2 = 6863558.206130736
$ws.Range("A18").Value2 = 111270755
$ws.Range("Q18").Value2 = 473194.7999623233
$ws.Range("R18").Value2 = 6863736.454484907
$ws.Range("A19").Value2 = 111271176
$ws.Range("B19").Value2 = 78579
$ws.Range("E19").Value2 = 2081
$ws.Range("F19").Value2 = "Skrovellav"
$ws.Range("G19").Value2 = "Lobaria scrobiculata"
$ws.Range("H19").Value2 = "(Scop.) DC."
$ws.Range("Q19").Value2 = 473227.9160841404
$ws.Range("R19").Value2 = 6863625.911539786
$ws.Range("A20").Value2 = 111271588
$ws.Range("B20").Value2 = 78578
$ws.Range("D20").Value2 = "NT"
$ws.Range("E20").Value2 = 6458
$ws.Range("F20").Value2 = "Lunglav"
$ws.Range("G20").Value2 = "Lobaria pulmonaria"
$ws.Range("H20").Value2 = "(L.) Hoffm."
$ws.Range("Q20").Value2 = 473140.3516782348
$ws.Range("R20").Value2 = 6863595.022241795
$ws.Range("A21").Value2 = 111271055
$ws.Range("B21").Value2 = 78579
$ws.Range("D21").Value2 = "NT"
$ws.Range("E21").Value2 = 2081
$ws.Range("F21").Value2 = "Skrovellav"
$ws.Range("G21").Value2 = "Lobaria scrobiculata"
$ws.Range("H21").Value2 = "(Scop.) DC."
$ws.Range("Q21").Value2 = 473238.8676645419
$ws.Range("R21").Value2 = 6863638.079474191
$ws.Range("A22").Value2 = 111271382
$ws.Range("B22").Value2 = 96348
$ws.Range("D22").Value2 = "VU"
$ws.Range("E22").Value2 = 220787
$ws.Range("F22").Value2 = "Knärot"
$ws.Range("G22").Value2 = "Goodyera repens"
$ws.Range("H22").Value2 = "(L.) R. Br."
$ws.Range("Q22").Value2 = 473167.6377000402
$ws.Range("R22").Value2 = 6863583.496200636
$ws.Range("AC22").Value2 = "Tre blommande."
$ws.Range("AC17").ClearContents()
$ws.Range("AC20").ClearContents()
$ws.Range("AC21").ClearContents()
